# el-254: Padronizacao da assinatura
#
# 1) Title paragraph: new spacing (before/after/line = 300, auto).
# 2) Remove the empty paragraph that used to sit right under the title.
# 3) Give the section a default header + footer (both start blank).
# 4) Tighten the top/bottom page margins now that the header/footer carry
#    some of that vertical space.

$d = $word.ActiveDocument

# --- 1) Title paragraph spacing -------------------------------------------
$title = $d.Paragraphs.Item(1)
$title.Format.SpaceBefore     = 15
$title.Format.SpaceAfter      = 15
$title.Format.LineSpacingRule = 0
$title.Format.LineSpacing     = 15

# --- 2) Drop the empty paragraph right after the title ---------------------
$blank = $d.Paragraphs.Item(2)
$blank.Range.Delete()

# --- 3) Header / footer -----------------------------------------------------
$section = $d.Sections.Item(1)

$header = $section.Headers.Item(1)
$header.LinkToPrevious = $false
$header.Range.Text = ""

$footer = $section.Footers.Item(1)
$footer.LinkToPrevious = $false
$footer.Range.Text = ""

# --- 4) Page margins ---------------------------------------------------------
$pageSetup = $section.PageSetup
$pageSetup.TopMargin    = 69.4488188976378
$pageSetup.BottomMargin = 45.35433070866142
